# eee3096_prac2.xlsx edit:
#  - Remove the "Attentuation" column (C) entirely, header + formulas.
#  - Append three more frequency-response measurement rows (20-22).
#  - Move the viewport/selection to the bottom of the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out column C (Attentuation header + ROUND() formulas) - the
# column formatting/width stays put, only the cell contents go away.
$ws.Range("C1:C19").ClearContents()

# New rows of measurements (OutputAmplitude, Frequency).
$newRows = @(
    @(1.7, 2000),
    @(1.5, 3000),
    @(1.1, 4000)
)

$row = 20
foreach ($pair in $newRows) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Update the view: scroll so row 11 is at the top, select the cell just
# past the new data (matches the author's saved selection state).
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("E22").Select() | Out-Null
